{"js": "// Insert a new \"List Bullet\" paragraph listing the two responsible\n// instructors right after the \"Docente(s) Respons\u00e1vel(eis)\" heading.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.trim() === \"Docente(s) Respons\u00e1vel(eis)\") {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Docente(s) Respons\u00e1vel(eis)' paragraph\");\n}\n\n// Build the new paragraph via a flat-OPC OOXML fragment so the run\n// structure matches exactly: first run carries the text followed by a\n// manual line break, second run carries the second name on its own.\nconst flatOpc = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"ListBullet\"/>\n            </w:pPr>\n            <w:r>\n              <w:t>3577649 - Carlos Angelo Nunes</w:t>\n              <w:br/>\n            </w:r>\n            <w:r>\n              <w:t>1922320 - Sebastiao Ribeiro</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nconst insertionPoint = target.getRange(\"End\");\ninsertionPoint.insertOoxml(flatOpc, \"After\");\n\nawait context.sync();\n", "ps1": "# Insert a new \"List Bullet\" paragraph listing the two responsible\n# instructors right after the \"Docente(s) Respons\u00e1vel(eis)\" heading.\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Trim() -eq \"Docente(s) Respons\u00e1vel(eis)\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not find the 'Docente(s) Respons\u00e1vel(eis)' paragraph\"\n}\n\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n\n# Build the new paragraph via a flat-OPC OOXML fragment so the run\n# structure matches exactly: first run carries the text followed by a\n# manual line break, second run carries the second name on its own.\n$flatOpc = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n    '<pkg:xmlData>' + `\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' + `\n    '<w:body>' + `\n    '<w:p>' + `\n    '<w:pPr><w:pStyle w:val=\"ListBullet\"/></w:pPr>' + `\n    '<w:r><w:t>3577649 - Carlos Angelo Nunes</w:t><w:br/></w:r>' + `\n    '<w:r><w:t>1922320 - Sebastiao Ribeiro</w:t></w:r>' + `\n    '</w:p>' + `\n    '</w:body>' + `\n    '</w:document>' + `\n    '</pkg:xmlData>' + `\n    '</pkg:part>' + `\n    '</pkg:package>'\n\n$newPara.Range.InsertXML($flatOpc)\n"}
